$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from I1 to J1 so the new header cell matches the existing bold/border style
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

# Apply cell value updates as described by the diff
$ws.Range("J1").Value = 'chunking_strategy'
$ws.Range("A2").Value = '2025-07-24 14:51:24'
$ws.Range("D2").Value = '國立中山大學新海研3號貴重儀器使用中心誠徵專任技術員1名'
$ws.Range("E2").Value = 'https://www.nstc.gov.tw/folksonomy/detail/ddc2e921-92c5-4004-8c2f-be2373c53f52?l=ch'
$ws.Range("F2").Value = '相關應徵資料予以保密，合者約談，不合者恕不另行通知。 發佈日期：2025-07-04 00:00:00'
$ws.Range("J2").Value = ''
$ws.Range("A3").Value = '2025-07-24 14:51:24'
$ws.Range("D3").Value = '[徵才] 國立臺灣大學防災減害與韌性學程 (綠‧韌性研究室) 徵求都市規劃/景觀/地理資訊專長 [專任計畫助理]'
$ws.Range("E3").Value = 'https://www.nstc.gov.tw/folksonomy/detail/2793c7ef-b68d-4f00-9388-e011b78b9553?l=ch'
$ws.Range("F3").Value = '3.其他有利申請之相關文件 發佈日期：2025-07-21 00:00:00'
$ws.Range("J3").Value = ''
$ws.Range("A4").Value = '2025-07-24 14:51:24'
$ws.Range("D4").Value = '中國醫藥大學 癌症生物精準醫學研究中心  王紹椿老師實驗室 誠徵 博士後研究員'
$ws.Range("E4").Value = 'https://www.nstc.gov.tw/folksonomy/detail/701ca4f1-a9f5-4a61-9b66-c4cf60f5c093?l=ch'
$ws.Range("F4").Value = '歡迎對癌症研究有興趣的夥伴加入我們的團隊！ 發佈日期：2025-07-14 00:00:00'
$ws.Range("J4").Value = ''
$ws.Range("A5").Value = '2025-07-24 14:51:24'
$ws.Range("D5").Value = '中國醫藥大學 癌症生物精準醫學研究中心  王紹椿老師實驗室 誠徵 碩士級研究助理'
$ws.Range("E5").Value = 'https://www.nstc.gov.tw/folksonomy/detail/2521ae27-55c0-4f27-9ded-b4bc908c1aff?l=ch'
$ws.Range("F5").Value = '歡迎對癌症研究有興趣的夥伴加入我們的團隊！ 發佈日期：2025-07-14 00:00:00'
$ws.Range("J5").Value = ''
$ws.Range("A6").Value = '2025-07-24 14:51:24'
$ws.Range("D6").Value = '國立臺東大學通識教育中心徵聘專任助理教授以上教師徵才公告，收件至114年8月15日止。'
$ws.Range("E6").Value = 'https://www.nstc.gov.tw/folksonomy/detail/e407fdbc-62c9-4e09-b08a-35a897cc4186?l=ch'
$ws.Range("F6").Value = '其    它： 相關訊息，請至本校首頁徵人啟事https://psn.nttu.edu.tw/p/406-1047-165359,r595.php?Lang=zh-tw查詢下載。 聯絡人姓名: 李家婕小姐 聯絡人電話: 089-517492 電子信箱：evalee@nttu.edu.tw 發佈日期：2025-07-09 00:00:00'
$ws.Range("J6").Value = ''
$ws.Range("A7").Value = '2025-07-24 15:35:07'
$ws.Range("D7").Value = '國立中山大學新海研3號貴重儀器使用中心誠徵專任技術員1名'
$ws.Range("E7").Value = 'https://www.nstc.gov.tw/folksonomy/detail/ddc2e921-92c5-4004-8c2f-be2373c53f52?l=ch'
$ws.Range("F7").Value = '相關應徵資料予以保密，合者約談，不合者恕不另行通知。 發佈日期：2025-07-04 00:00:00'
$ws.Range("J7").Value = 'hybrid_chunking'
$ws.Range("A8").Value = '2025-07-24 15:35:07'
$ws.Range("D8").Value = '[徵才] 國立臺灣大學防災減害與韌性學程 (綠‧韌性研究室) 徵求都市規劃/景觀/地理資訊專長 [專任計畫助理]'
$ws.Range("E8").Value = 'https://www.nstc.gov.tw/folksonomy/detail/2793c7ef-b68d-4f00-9388-e011b78b9553?l=ch'
$ws.Range("F8").Value = '3.其他有利申請之相關文件 發佈日期：2025-07-21 00:00:00'
$ws.Range("J8").Value = 'hybrid_chunking'
$ws.Range("A9").Value = '2025-07-24 15:35:07'
$ws.Range("D9").Value = '中國醫藥大學 癌症生物精準醫學研究中心  王紹椿老師實驗室 誠徵 博士後研究員'
$ws.Range("E9").Value = 'https://www.nstc.gov.tw/folksonomy/detail/701ca4f1-a9f5-4a61-9b66-c4cf60f5c093?l=ch'
$ws.Range("F9").Value = '歡迎對癌症研究有興趣的夥伴加入我們的團隊！ 發佈日期：2025-07-14 00:00:00'
$ws.Range("J9").Value = 'hybrid_chunking'
$ws.Range("A10").Value = '2025-07-24 15:35:07'
$ws.Range("D10").Value = '中國醫藥大學 癌症生物精準醫學研究中心  王紹椿老師實驗室 誠徵 碩士級研究助理'
$ws.Range("E10").Value = 'https://www.nstc.gov.tw/folksonomy/detail/2521ae27-55c0-4f27-9ded-b4bc908c1aff?l=ch'
$ws.Range("F10").Value = '歡迎對癌症研究有興趣的夥伴加入我們的團隊！ 發佈日期：2025-07-14 00:00:00'
$ws.Range("J10").Value = 'hybrid_chunking'
$ws.Range("A11").Value = '2025-07-24 15:35:07'
$ws.Range("D11").Value = '國立臺東大學通識教育中心徵聘專任助理教授以上教師徵才公告，收件至114年8月15日止。'
$ws.Range("E11").Value = 'https://www.nstc.gov.tw/folksonomy/detail/e407fdbc-62c9-4e09-b08a-35a897cc4186?l=ch'
$ws.Range("F11").Value = '其    它： 相關訊息，請至本校首頁徵人啟事https://psn.nttu.edu.tw/p/406-1047-165359,r595.php?Lang=zh-tw查詢下載。 聯絡人姓名: 李家婕小姐 聯絡人電話: 089-517492 電子信箱：evalee@nttu.edu.tw 發佈日期：2025-07-09 00:00:00'
$ws.Range("J11").Value = 'hybrid_chunking'
$ws.Range("A12").Value = '2025-07-24 15:35:39'
$ws.Range("B12").Value = '材料相關的職缺有哪些？'
$ws.Range("D12").Value = '國立中山大學新海研3號貴重儀器使用中心誠徵專任技術員1名'
$ws.Range("E12").Value = 'https://www.nstc.gov.tw/folksonomy/detail/ddc2e921-92c5-4004-8c2f-be2373c53f52?l=ch'
$ws.Range("F12").Value = '相關應徵資料予以保密，合者約談，不合者恕不另行通知。 發佈日期：2025-07-04 00:00:00'
$ws.Range("J12").Value = 'hybrid_chunking'
$ws.Range("A13").Value = '2025-07-24 15:35:39'
$ws.Range("B13").Value = '材料相關的職缺有哪些？'
$ws.Range("D13").Value = '[徵才] 國立臺灣大學防災減害與韌性學程 (綠‧韌性研究室) 徵求都市規劃/景觀/地理資訊專長 [專任計畫助理]'
$ws.Range("E13").Value = 'https://www.nstc.gov.tw/folksonomy/detail/2793c7ef-b68d-4f00-9388-e011b78b9553?l=ch'
$ws.Range("F13").Value = '3.其他有利申請之相關文件 發佈日期：2025-07-21 00:00:00'
$ws.Range("J13").Value = 'hybrid_chunking'
$ws.Range("A14").Value = '2025-07-24 15:35:39'
$ws.Range("B14").Value = '材料相關的職缺有哪些？'
$ws.Range("D14").Value = '中國醫藥大學 癌症生物精準醫學研究中心  王紹椿老師實驗室 誠徵 博士後研究員'
$ws.Range("E14").Value = 'https://www.nstc.gov.tw/folksonomy/detail/701ca4f1-a9f5-4a61-9b66-c4cf60f5c093?l=ch'
$ws.Range("F14").Value = '歡迎對癌症研究有興趣的夥伴加入我們的團隊！ 發佈日期：2025-07-14 00:00:00'
$ws.Range("J14").Value = 'hybrid_chunking'
$ws.Range("A15").Value = '2025-07-24 15:35:39'
$ws.Range("B15").Value = '材料相關的職缺有哪些？'
$ws.Range("D15").Value = '中國醫藥大學 癌症生物精準醫學研究中心  王紹椿老師實驗室 誠徵 碩士級研究助理'
$ws.Range("E15").Value = 'https://www.nstc.gov.tw/folksonomy/detail/2521ae27-55c0-4f27-9ded-b4bc908c1aff?l=ch'
$ws.Range("F15").Value = '歡迎對癌症研究有興趣的夥伴加入我們的團隊！ 發佈日期：2025-07-14 00:00:00'
$ws.Range("J15").Value = 'hybrid_chunking'
$ws.Range("A16").Value = '2025-07-24 15:35:39'
$ws.Range("B16").Value = '材料相關的職缺有哪些？'
$ws.Range("D16").Value = '國立臺東大學通識教育中心徵聘專任助理教授以上教師徵才公告，收件至114年8月15日止。'
$ws.Range("E16").Value = 'https://www.nstc.gov.tw/folksonomy/detail/e407fdbc-62c9-4e09-b08a-35a897cc4186?l=ch'
$ws.Range("F16").Value = '其    它： 相關訊息，請至本校首頁徵人啟事https://psn.nttu.edu.tw/p/406-1047-165359,r595.php?Lang=zh-tw查詢下載。 聯絡人姓名: 李家婕小姐 聯絡人電話: 089-517492 電子信箱：evalee@nttu.edu.tw 發佈日期：2025-07-09 00:00:00'
$ws.Range("J16").Value = 'hybrid_chunking'

Write-Output "done"
